# Apply update described by the diff:
# - Rows 2 and 4 swap their "species" identity data (A, E, F, G, H) and get new
#   B (Taxonsorteringsordning) values.
# - Row 3 gets a new B value only.
# - Four brand-new rows (5-8) are appended with new observation records.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2 (now holds the former row-4 "Dofttaggsvamp" species data) ----
$ws.Range("A2").Value = 112206837
$ws.Range("B2").Value = 90818
$ws.Range("E2").Value = 4368
$ws.Range("F2").Value = "Dofttaggsvamp"
$ws.Range("G2").Value = "Hydnellum suaveolens"
$ws.Range("H2").Value = "(Scop.:Fr.) P. Karst."

# ---- Row 3 (only the Taxonsorteringsordning changes) ----
$ws.Range("B3").Value = 90802

# ---- Row 4 (now holds the former row-2 "Orange taggsvamp" species data) ----
$ws.Range("A4").Value = 112206846
$ws.Range("B4").Value = 90792
$ws.Range("E4").Value = 4361
$ws.Range("F4").Value = "Orange taggsvamp"
$ws.Range("G4").Value = "Hydnellum aurantiacum"
$ws.Range("H4").Value = "(Batsch:Fr.) P.Karst."

# ---- New row 5 ----
$ws.Range("A5").Value = 112395257
$ws.Range("B5").Value = 90802
$ws.Range("C5").Value = "Ovaliderad"
$ws.Range("D5").Value = "NT"
$ws.Range("E5").Value = 788
$ws.Range("F5").Value = "Gul taggsvamp"
$ws.Range("G5").Value = "Hydnellum geogenium"
$ws.Range("H5").Value = "(Fr.) Banker"
$ws.Range("P5").Value = "Påterud, Vrm"
$ws.Range("Q5").Value = 333022
$ws.Range("R5").Value = 6626625
$ws.Range("S5").Value = 10
$ws.Range("T5").Value = "Värmland"
$ws.Range("U5").Value = "Eda"
$ws.Range("V5").Value = "Värmland"
$ws.Range("W5").Value = "Järnskog"
$ws.Range("Y5").Value = "'2023-09-19"
$ws.Range("AA5").Value = "'2023-09-19"
$ws.Range("AD5").Value = $false
$ws.Range("AE5").Value = $false
$ws.Range("AG5").Value = $false
$ws.Range("AW5").Value = "Jan Rees"
$ws.Range("AX5").Value = "Jan Rees"

# ---- New row 6 ----
$ws.Range("A6").Value = 112395259
$ws.Range("B6").Value = 90792
$ws.Range("C6").Value = "Ovaliderad"
$ws.Range("D6").Value = "NT"
$ws.Range("E6").Value = 4361
$ws.Range("F6").Value = "Orange taggsvamp"
$ws.Range("G6").Value = "Hydnellum aurantiacum"
$ws.Range("H6").Value = "(Batsch:Fr.) P.Karst."
$ws.Range("P6").Value = "Påterud, Vrm"
$ws.Range("Q6").Value = 333042
$ws.Range("R6").Value = 6626628
$ws.Range("S6").Value = 10
$ws.Range("T6").Value = "Värmland"
$ws.Range("U6").Value = "Eda"
$ws.Range("V6").Value = "Värmland"
$ws.Range("W6").Value = "Järnskog"
$ws.Range("Y6").Value = "'2023-09-19"
$ws.Range("AA6").Value = "'2023-09-19"
$ws.Range("AD6").Value = $false
$ws.Range("AE6").Value = $false
$ws.Range("AG6").Value = $false
$ws.Range("AW6").Value = "Jan Rees"
$ws.Range("AX6").Value = "Jan Rees"

# ---- New row 7 ----
$ws.Range("A7").Value = 112395256
$ws.Range("B7").Value = 93293
$ws.Range("C7").Value = "Ovaliderad"
$ws.Range("D7").Value = "LC"
$ws.Range("E7").Value = 2671
$ws.Range("F7").Value = "Fällmossa"
$ws.Range("G7").Value = "Antitrichia curtipendula"
$ws.Range("H7").Value = "(Hedw.) Brid."
$ws.Range("P7").Value = "Påterud, Vrm"
$ws.Range("Q7").Value = 333013
$ws.Range("R7").Value = 6626631
$ws.Range("S7").Value = 10
$ws.Range("T7").Value = "Värmland"
$ws.Range("U7").Value = "Eda"
$ws.Range("V7").Value = "Värmland"
$ws.Range("W7").Value = "Järnskog"
$ws.Range("Y7").Value = "'2023-09-19"
$ws.Range("AA7").Value = "'2023-09-19"
$ws.Range("AD7").Value = $false
$ws.Range("AE7").Value = $false
$ws.Range("AG7").Value = $false
$ws.Range("AW7").Value = "Jan Rees"
$ws.Range("AX7").Value = "Jan Rees"

# ---- New row 8 ----
$ws.Range("A8").Value = 112395260
$ws.Range("B8").Value = 90818
$ws.Range("C8").Value = "Ovaliderad"
$ws.Range("D8").Value = "NT"
$ws.Range("E8").Value = 4368
$ws.Range("F8").Value = "Dofttaggsvamp"
$ws.Range("G8").Value = "Hydnellum suaveolens"
$ws.Range("H8").Value = "(Scop.:Fr.) P. Karst."
$ws.Range("P8").Value = "Påterud, Vrm"
$ws.Range("Q8").Value = 333038
$ws.Range("R8").Value = 6626631
$ws.Range("S8").Value = 10
$ws.Range("T8").Value = "Värmland"
$ws.Range("U8").Value = "Eda"
$ws.Range("V8").Value = "Värmland"
$ws.Range("W8").Value = "Järnskog"
$ws.Range("Y8").Value = "'2023-09-19"
$ws.Range("AA8").Value = "'2023-09-19"
$ws.Range("AD8").Value = $false
$ws.Range("AE8").Value = $false
$ws.Range("AG8").Value = $false
$ws.Range("AW8").Value = "Jan Rees"
$ws.Range("AX8").Value = "Jan Rees"
